$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 16 (Card ID 14): "Eldritch Horror" changes from MINION (999 HP) to LORD ---
$ws.Range("J16").Value = 'If you control 3 occult devotees at the end of your turn, summon this card with 30 health and 30 strength.'
$ws.Range("F16").Value = "LORD"
$ws.Range("G16").Value = 0

# --- Row 15 (Card ID 13): "Occult Devotee" becomes "Eldritch Devotee" ---
$ws.Range("E15").Value = "Eldritch Devotee"
$ws.Range("J15").Value = 'At the end of your turn, search 1 "Devoted Cultist"'

# --- Row 40 (Card ID 38): new "Assasination" card ---
$ws.Range("E40").Value = "Assasination"
$ws.Range("J40").Value = "Kill an undamaged unit."
$ws.Range("F40").Value = "UTILITY"
$ws.Range("G40").Value = 5
$ws.Range("H40").Value = 0
$ws.Range("I40").Value = 0
$ws.Range("K40").Value = "N"

# --- Row 38 (Card ID 36): "The Healer" renamed to "The Doctor" ---
$ws.Range("E38").Value = "The Doctor"

# --- Update the saved view state (top-left cell and active selection) ---
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 4
$ws.Range("E38").Select()
